$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.741.03'
$ws.Range("E2").Value = '  -1.56%  '
$ws.Range("D3").Value = '2.445.95'
$ws.Range("E3").Value = '  -3.01%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '522.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("D9").Value = '2.448.55'
$ws.Range("E9").Value = '  -2.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0972'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("E11").Value = '  -3.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.90'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.16%  '
$ws.Range("E13").Value = '  -3.59%  '
$ws.Range("D14").Value = '2.879.36'
$ws.Range("E14").Value = '  -2.01%  '
$ws.Range("D15").Value = '57.652.05'
$ws.Range("E15").Value = '  -1.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.65'
$ws.Range("D16").Style = "Normal"
$ws.Range("E17").Value = '  -2.03%  '
$ws.Range("D18").Value = '2.445.28'
$ws.Range("E18").Value = '  -2.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.13%  '
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '314.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.58%  '
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.81%  '
$ws.Range("E25").Value = '  +1.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.47%  '
$ws.Range("E27").Value = '  -3.08%  '
$ws.Range("E28").Value = '  -2.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '172.50'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.75%  '
$ws.Range("D30").Value = '0.0₃0733'
$ws.Range("E30").Value = '  -2.90%  '
$ws.Range("E31").Value = '  -1.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.78%  '
$ws.Range("E33").Value = '  -5.51%  '
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.80'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.57%  '
$ws.Range("E37").Value = '  -6.41%  '
$ws.Range("E38").Value = '  -4.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.28'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("E40").Value = '  -1.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.787'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.40'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '263.62'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.583'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.80'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0927'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.66%  '
$ws.Range("E48").Value = '  -2.07%  '
$ws.Range("E49").Value = '  -1.82%  '
$ws.Range("E50").Value = '  -4.45%  '
